$wb = $excel.ActiveWorkbook

function Set-TextValue($Range, $Text) {
    # Force the cell to be stored as a shared string (text) rather than a
    # number, while preserving the original "General"/vertical-top cell
    # style (no new style record should be created).
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
    $Range.VerticalAlignment = -4160   # xlTop
}

$sheetNames = @("Sheet1", "Sheet2", "Sheet3", "Sheet4")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-TextValue $ws.Range("F2")  "9840067290"
    Set-TextValue $ws.Range("AE2") "9840054005"
    Set-TextValue $ws.Range("AT2") "9840083579"
    Set-TextValue $ws.Range("AX2") "9840060672"
}
